$d = $word.ActiveDocument

# 1. Global font replace: TimesNewToman -> Times New Roman across entire document
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# 2. Title paragraph
$p = $d.Paragraphs.Item(1)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = 'The Beauty and Wonders of Chemistry: Unveiling the Molecular Symphony'

# 3. Author paragraph
$p = $d.Paragraphs.Item(2)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = 'Dr. Alice Thompson'

# 4. Email paragraph
$p = $d.Paragraphs.Item(3)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = 'thompson.alice@schoolmail.edu'

# paragraph 4 is already empty, nothing to do

# 5. Main body paragraph with line breaks
$p = $d.Paragraphs.Item(5)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = 'Chemistry, the study of matter and its transformations, is a fascinating field that unveils the intricate molecular symphony that orchestrates the world around us. It offers a profound understanding of how substances interact, empowering us to harness their properties and engineer materials with remarkable applications. Chemistry''s reach extends from the macroscopic realm of everyday phenomena to the enigmatic quantum realm, where particles exhibit behaviors that defy classical intuition.Chemistry unravels the mysteries of life, illuminating the intricate biochemical pathways that govern cellular processes. It unravels the complexities of materials, revealing the atomic structures and interactions that determine their properties. This knowledge fuels technological advancements, leading to innovative materials, medicines, and energy sources that shape our modern world. Delving into chemistry is embarking on an intellectual journey that unveils the elegance and wonder of the molecular world.The study of chemistry fosters critical thinking, analytical reasoning, and problem-solving abilities, equipping students with valuable skills that extend beyond the laboratory. It encourages a spirit of inquiry, curiosity, and exploration, nurturing the scientific mindset that drives progress and innovation. Whether aspiring to careers in science, medicine, engineering, or beyond, a solid foundation in chemistry provides a gateway to countless opportunities.'

# 6. Summary heading paragraph (unchanged text, but still has TimesNewToman which is already handled globally)

# 7. Summary body paragraph; also drops lastRenderedPageBreak automatically
$p = $d.Paragraphs.Item(7)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = 'Chemistry, the study of matter and its transformations, unveils the intricacies of the molecular world. It offers a profound understanding of substances, their interactions, and their applications, shaping our understanding of life, materials, and technology. The study of chemistry not only imparts knowledge but also cultivates critical thinking, analytical reasoning, and problem-solving abilities, empowering students to navigate an increasingly complex and scientifically driven world. Embarking on this intellectual adventure reveals the elegance and wonder of the molecular symphony that orchestrates the universe.'

# 8. Append new empty trailing paragraph
$d.Content.InsertParagraphAfter()

Write-Output "Edit complete"
